$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns (D, E) sometimes contain values that look like
# plain numbers (e.g. "0.4823") even though the column stores text (e.g.
# "31.156.03" is not a valid number). Force the cells to be treated as text
# while they are populated, then restore the original (default) style so the
# cell formatting is left exactly as it was.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "31.156.03"
$ws.Range("E2").Value = "  +1.28%  "

$ws.Range("D3").Value = "1.944.54"
$ws.Range("E3").Value = "  +0.87%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "242.53"

$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").Value = "0.4823"
$ws.Range("E7").Value = "  -0.50%  "

$ws.Range("D8").Value = "0.2914"
$ws.Range("E8").Value = "  -0.55%  "

$ws.Range("D9").Value = "0.06805"

$ws.Range("D10").Value = "20.08"
$ws.Range("E10").Value = "  +4.95%  "

$ws.Range("D11").Value = "104.50"
$ws.Range("E11").Value = "  -1.61%  "

$ws.Range("D12").Value = "0.07855"
$ws.Range("E12").Value = "  +1.11%  "

$ws.Range("D13").Value = "1.957.50"
$ws.Range("E13").Value = "  +1.52%  "

$ws.Range("D14").Value = "5.315"
$ws.Range("E14").Value = "  -0.03%  "

$ws.Range("D15").Value = "0.6907"
$ws.Range("E15").Value = "  -0.71%  "

$ws.Range("D16").Value = "296.27"
$ws.Range("E16").Value = "  +7.75%  "

$ws.Range("D17").Value = "31.170.57"
$ws.Range("E17").Value = "  +1.32%  "

$ws.Range("D18").Value = "2.212.41"
$ws.Range("E18").Value = "  +1.20%  "

$ws.Range("D19").Value = "13.01"
$ws.Range("E19").Value = "  +0.64%  "

$ws.Range("D20").Value = "0.000007619"
$ws.Range("E20").Value = "  -0.48%  "

$ws.Range("D21").Value = "5.589"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "6.466"
$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "9.582"
$ws.Range("E25").Value = "  -2.59%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "169.21"
$ws.Range("E26").Value = "  +2.88%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "19.85"
$ws.Range("E27").Value = "  +2.22%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.143"
$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "1.396"
$ws.Range("E29").Value = "  +1.02%  "

$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "0.1017"
$ws.Range("E30").Value = "  -1.90%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "4.629"
$ws.Range("E31").Value = "  +1.53%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "1.537"
$ws.Range("E32").Value = "  -0.66%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "4.365"
$ws.Range("E33").Value = "  +0.48%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.04845"
$ws.Range("E34").Value = "  -0.49%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.7432"
$ws.Range("E35").Value = "  -1.61%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.134"
$ws.Range("E36").Value = "  -0.54%  "

$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.745"
$ws.Range("E37").Value = "  +1.02%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01962"
$ws.Range("E38").Value = "  -1.12%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "6.612"
$ws.Range("E39").Value = "  +2.69%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.655"
$ws.Range("E40").Value = "  +0.36%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "77.46"
$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "2.034"
$ws.Range("E42").Value = "  -0.89%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.8756"
$ws.Range("E43").Value = "  -0.61%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.4373"
$ws.Range("E44").Value = "  -0.98%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "106.57"
$ws.Range("E45").Value = "  -0.80%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.020.18"
$ws.Range("E47").Value = "  +4.03%  "

$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "7.575"
$ws.Range("E48").Value = "  -3.35%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.229"
$ws.Range("E49").Value = "  +0.42%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.1215"
$ws.Range("E50").Value = "  -1.52%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "35.21"
$ws.Range("E51").Value = "  -2.29%  "

# Restore original (default, unstyled) formatting for the data range so the
# cell style metadata matches the source workbook.
$dataRange.Style = "Normal"
